# Generate Report for Handback
#
# For each localized-language sheet (zh-cn, de-de) this marks the two
# tracked files as handed back: the Status column flips from
# "Ready for handoff" to "Handed back: in sync with en-US", the
# "Latest Target File" / "Latest Handback File" columns get populated
# (mirroring the source markdown file and the handoff .xlf that was
# round-tripped back), and "Latest Handback DateTime" is stamped with
# the handback time. The Overview sheet shares the same Status text so
# it picks up the new wording automatically.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------------
# B2/C2/B3/C3 just mirror the per-language Status values.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusHandedBack
$overview.Range("C2").Value = $statusHandedBack
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet -----------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("B2").Value = $statusHandedBack
$zhcn.Range("E2").Value = "e036e640-2d81-4bf6-a212-a69565d14694.md"
$zhcn.Range("E2").Style = "HyperLink"
$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/441c9c2c4ddf032da7cbe32aad436f04c44f1917/e2e/e036e640-2d81-4bf6-a212-a69565d14694.md", "", "", "e036e640-2d81-4bf6-a212-a69565d14694.md")
$zhcn.Range("F2").Value = "e036e640-2d81-4bf6-a212-a69565d14694.16733f7d71ef6f06c9a1a7720778649cd2e1b81b.zh-cn.xlf"
$zhcn.Range("F2").Style = "HyperLink"
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1ababef8e7e397e40e4cd4d0118d9e1ef18cd88b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/e036e640-2d81-4bf6-a212-a69565d14694.16733f7d71ef6f06c9a1a7720778649cd2e1b81b.zh-cn.xlf", "", "", "e036e640-2d81-4bf6-a212-a69565d14694.16733f7d71ef6f06c9a1a7720778649cd2e1b81b.zh-cn.xlf")
$zhcn.Range("G2").Value = "2016-03-03 09:06:18"

$zhcn.Range("B3").Value = $statusHandedBack
$zhcn.Range("E3").Value = "333085f4-f54c-4a78-8c04-096cc915e0fe.md"
$zhcn.Range("E3").Style = "HyperLink"
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/441c9c2c4ddf032da7cbe32aad436f04c44f1917/e2e/333085f4-f54c-4a78-8c04-096cc915e0fe.md", "", "", "333085f4-f54c-4a78-8c04-096cc915e0fe.md")
$zhcn.Range("F3").Value = "333085f4-f54c-4a78-8c04-096cc915e0fe.2e78b0eb2414c3bd9be7fdb837a269b99732dcf7.zh-cn.xlf"
$zhcn.Range("F3").Style = "HyperLink"
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1ababef8e7e397e40e4cd4d0118d9e1ef18cd88b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/333085f4-f54c-4a78-8c04-096cc915e0fe.2e78b0eb2414c3bd9be7fdb837a269b99732dcf7.zh-cn.xlf", "", "", "333085f4-f54c-4a78-8c04-096cc915e0fe.2e78b0eb2414c3bd9be7fdb837a269b99732dcf7.zh-cn.xlf")
$zhcn.Range("G3").Value = "2016-03-03 09:06:18"

# --- de-de sheet -----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("B2").Value = $statusHandedBack
$dede.Range("E2").Value = "e036e640-2d81-4bf6-a212-a69565d14694.md"
$dede.Range("E2").Style = "HyperLink"
$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/441c9c2c4ddf032da7cbe32aad436f04c44f1917/e2e/e036e640-2d81-4bf6-a212-a69565d14694.md", "", "", "e036e640-2d81-4bf6-a212-a69565d14694.md")
$dede.Range("F2").Value = "e036e640-2d81-4bf6-a212-a69565d14694.16733f7d71ef6f06c9a1a7720778649cd2e1b81b.de-de.xlf"
$dede.Range("F2").Style = "HyperLink"
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5fdc16e768e48208bbe7a1c4d90d4e93e3223a3a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/e036e640-2d81-4bf6-a212-a69565d14694.16733f7d71ef6f06c9a1a7720778649cd2e1b81b.de-de.xlf", "", "", "e036e640-2d81-4bf6-a212-a69565d14694.16733f7d71ef6f06c9a1a7720778649cd2e1b81b.de-de.xlf")
$dede.Range("G2").Value = "2016-03-03 09:06:39"

$dede.Range("B3").Value = $statusHandedBack
$dede.Range("E3").Value = "333085f4-f54c-4a78-8c04-096cc915e0fe.md"
$dede.Range("E3").Style = "HyperLink"
$dede.Hyperlinks.Add($dede.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/441c9c2c4ddf032da7cbe32aad436f04c44f1917/e2e/333085f4-f54c-4a78-8c04-096cc915e0fe.md", "", "", "333085f4-f54c-4a78-8c04-096cc915e0fe.md")
$dede.Range("F3").Value = "333085f4-f54c-4a78-8c04-096cc915e0fe.2e78b0eb2414c3bd9be7fdb837a269b99732dcf7.de-de.xlf"
$dede.Range("F3").Style = "HyperLink"
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5fdc16e768e48208bbe7a1c4d90d4e93e3223a3a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/333085f4-f54c-4a78-8c04-096cc915e0fe.2e78b0eb2414c3bd9be7fdb837a269b99732dcf7.de-de.xlf", "", "", "333085f4-f54c-4a78-8c04-096cc915e0fe.2e78b0eb2414c3bd9be7fdb837a269b99732dcf7.de-de.xlf")
$dede.Range("G3").Value = "2016-03-03 09:06:39"
